$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.835.95'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.31%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.567.71'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.14%  '

# Row 4
$ws.Range("E4").Value = '  -0.02%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '206.73'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.91%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.491'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.73%  '

# Row 7
$ws.Range("E7").Value = '  +0.03%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.00'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.86%  '

# Row 9
$ws.Range("E9").Value = '  -0.75%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0587'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.94%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0864'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.16%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.789.69'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.22%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.569.81'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.02%  '

# Row 14
$ws.Range("E14").Value = '  -2.21%  '

# Row 15
$ws.Range("E15").Value = '  -0.67%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '26.841.72'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.31%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.53'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.56%  '

# Row 18
$ws.Range("E18").Value = '  +1.76%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '214.71'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.18%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0678'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.94%  '

# Row 21
$ws.Range("E21").Value = '  +0.03%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.13'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.16%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.30'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.72%  '

# Row 24
$ws.Range("E24").Value = '  -1.01%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.60'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.64%  '

# Row 26
$ws.Range("E26").Value = '  +0.09%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '14.99'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.27%  '

# Row 28
$ws.Range("E28").Value = '  +0.01%  '

# Row 29
$ws.Range("E29").Value = '  -1.10%  '

# Row 30
$ws.Range("E30").Value = '  -0.69%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.12'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.27%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.18'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.79%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.403.34'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.61%  '

# Row 34
$ws.Range("E34").Value = '  -1.69%  '

# Row 35
$ws.Range("E35").Value = '  -0.96%  '

# Row 36
$ws.Range("E36").Value = '  -1.19%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.937'
$ws.Range("D37").Style = "Normal"

# Row 38
$ws.Range("E38").Value = '  -2.82%  '

# Row 39
$ws.Range("E39").Value = '  -2.77%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.815'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.53%  '

# Row 42
$ws.Range("E42").Value = '  +1.09%  '

# Row 43
$ws.Range("E43").Value = '  -0.48%  '

# Row 44
$ws.Range("B44").Value = 'MXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.18'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.96%  '

# Row 45
$ws.Range("B45").Value = 'FraxShare'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.31'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.93%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '63.36'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.35%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.703.29'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.01%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '86.18'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.79%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0₇0982'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.04%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0955'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.14%  '

# Row 51
$ws.Range("E51").Value = '  -0.98%  '
